# Update cryptos list with latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '29.207.48'
Set-TextCell 2 5 '  -0.47%  '

Set-TextCell 3 4 '1.859.38'
Set-TextCell 3 5 '  -0.82%  '

Set-TextCell 4 4 '0.9997'
Set-TextCell 4 5 '  -0.54%  '

Set-TextCell 5 4 '0.7153'
Set-TextCell 5 5 '  +0.84%  '

Set-TextCell 6 4 '241.44'
Set-TextCell 6 5 '  -0.78%  '

Set-TextCell 7 4 '0.9996'
Set-TextCell 7 5 '  -0.66%  '

Set-TextCell 8 4 '0.07808'
Set-TextCell 8 5 '  -0.42%  '

Set-TextCell 9 4 '0.3110'
Set-TextCell 9 5 '  -0.04%  '

Set-TextCell 10 5 '  -2.05%  '

Set-TextCell 11 4 '0.07811'
Set-TextCell 11 5 '  -3.35%  '

Set-TextCell 12 4 '1.850.03'
Set-TextCell 12 5 '  -1.69%  '

Set-TextCell 13 2 'Polkadot'
Set-TextCell 13 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 13 4 '5.115'
Set-TextCell 13 5 '  -1.05%  '

Set-TextCell 14 2 'Litecoin'
Set-TextCell 14 3 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 14 4 '92.33'
Set-TextCell 14 5 '  -1.19%  '

Set-TextCell 15 4 '0.6889'
Set-TextCell 15 5 '  -1.49%  '

Set-TextCell 16 4 '6.527'
Set-TextCell 16 5 '  +2.55%  '

Set-TextCell 17 4 '0.000008466'
Set-TextCell 17 5 '  +1.74%  '

Set-TextCell 18 4 '29.195.15'

Set-TextCell 19 4 '250.49'
Set-TextCell 19 5 '  -0.33%  '

Set-TextCell 20 4 '2.097.18'
Set-TextCell 20 5 '  -2.16%  '

Set-TextCell 21 4 '12.88'
Set-TextCell 21 5 '  -2.51%  '

Set-TextCell 22 5 '  -0.40%  '

Set-TextCell 23 4 '7.521'
Set-TextCell 23 5 '  -0.90%  '

Set-TextCell 24 4 '0.9998'
Set-TextCell 24 5 '  -0.60%  '

Set-TextCell 25 4 '0.1556'
Set-TextCell 25 5 '  -0.65%  '

Set-TextCell 26 4 '160.21'
Set-TextCell 26 5 '  -0.62%  '

Set-TextCell 27 4 '8.878'
Set-TextCell 27 5 '  -1.27%  '

Set-TextCell 28 4 '18.56'
Set-TextCell 28 5 '  -0.99%  '

Set-TextCell 29 4 '1.563'
Set-TextCell 29 5 '  +3.50%  '

Set-TextCell 30 4 '4.273'
Set-TextCell 30 5 '  -1.43%  '

Set-TextCell 31 5 '  -0.63%  '

Set-TextCell 32 4 '1.206'
Set-TextCell 32 5 '  -1.94%  '

Set-TextCell 33 4 '0.05210'
Set-TextCell 33 5 '  -1.04%  '

Set-TextCell 34 4 '0.7596'
Set-TextCell 34 5 '  +1.35%  '

Set-TextCell 35 4 '1.173'
Set-TextCell 35 5 '  +0.48%  '

Set-TextCell 36 4 '1.847'
Set-TextCell 36 5 '  -3.16%  '

Set-TextCell 37 4 '2.710'
Set-TextCell 37 5 '  -0.18%  '

Set-TextCell 38 4 '0.01860'
Set-TextCell 38 5 '  -0.31%  '

Set-TextCell 39 4 '1.226.52'
Set-TextCell 39 5 '  -3.24%  '

Set-TextCell 40 4 '2.730'
Set-TextCell 40 5 '  -1.04%  '

Set-TextCell 41 4 '0.8991'
Set-TextCell 41 5 '  -0.13%  '

Set-TextCell 42 4 '109.39'
Set-TextCell 42 5 '  -1.71%  '

Set-TextCell 43 4 '0.9989'
Set-TextCell 43 5 '  -0.77%  '

Set-TextCell 44 4 '5.648'
Set-TextCell 44 5 '  -10.24%  '

Set-TextCell 45 4 '1.996.97'
Set-TextCell 45 5 '  -1.88%  '

Set-TextCell 46 5 '  -0.83%  '

Set-TextCell 47 2 'Mantle'
Set-TextCell 47 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 47 4 '0.5183'
Set-TextCell 47 5 '  -0.76%  '

Set-TextCell 48 2 'EnergySwap'
Set-TextCell 48 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 48 4 '9.535'
Set-TextCell 48 5 '  +2.01%  '

Set-TextCell 49 2 'Aave'
Set-TextCell 49 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 49 4 '64.67'
Set-TextCell 49 5 '  -9.70%  '

Set-TextCell 50 4 '1.756'
Set-TextCell 50 5 '  -1.73%  '

Set-TextCell 51 4 '7.028'
Set-TextCell 51 5 '  +0.24%  '
